$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format cells whose new price values look numeric so Excel keeps them as literal text
# (matches the source data which stores these as plain strings, not numbers).
$textCells = @("D4", "D6", "D7", "D8", "D9", "D10", "D13", "D14", "D15", "D16", "D17", "D18", "D20", "D23", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "26.388.07"
$ws.Range("E2").Value = "  -2.62%  "

# Row 3
$ws.Range("D3").Value = "1.774.49"
$ws.Range("E3").Value = "  -1.32%  "

# Row 4
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.72%  "

# Row 5
$ws.Range("E5").Value = "  -0.64%  "

# Row 6
$ws.Range("D6").Value = "306.71"
$ws.Range("E6").Value = "  -0.44%  "

# Row 7
$ws.Range("D7").Value = "0.4237"
$ws.Range("E7").Value = "  +1.78%  "

# Row 8
$ws.Range("D8").Value = "0.3604"
$ws.Range("E8").Value = "  +1.54%  "

# Row 9
$ws.Range("D9").Value = "0.07147"
$ws.Range("E9").Value = "  +1.74%  "

# Row 10
$ws.Range("D10").Value = "0.8376"
$ws.Range("E10").Value = "  -0.42%  "

# Row 11
$ws.Range("E11").Value = "  +1.37%  "

# Row 12
$ws.Range("D12").Value = "1.763.00"
$ws.Range("E12").Value = "  -9.06%  "

# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "5.249"
$ws.Range("E13").Value = "  -0.10%  "

# Row 14
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "6.439"
$ws.Range("E14").Value = "  +1.87%  "

# Row 15
$ws.Range("D15").Value = "0.06894"
$ws.Range("E15").Value = "  +1.10%  "

# Row 16
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  -0.97%  "

# Row 17
$ws.Range("D17").Value = "79.02"
$ws.Range("E17").Value = "  -0.92%  "

# Row 18
$ws.Range("D18").Value = "0.000008659"
$ws.Range("E18").Value = "  -0.31%  "

# Row 19
$ws.Range("E19").Value = "  -0.84%  "

# Row 20
$ws.Range("D20").Value = "14.89"
$ws.Range("E20").Value = "  -1.03%  "

# Row 21
$ws.Range("D21").Value = "26.400.57"
$ws.Range("E21").Value = "  -4.56%  "

# Row 22
$ws.Range("E22").Value = "  +0.97%  "

# Row 23
$ws.Range("D23").Value = "10.90"
$ws.Range("E23").Value = "  +1.67%  "

# Row 24
$ws.Range("D24").Value = "1.981.49"
$ws.Range("E24").Value = "  -7.59%  "

# Row 25
$ws.Range("D25").Value = "151.48"
$ws.Range("E25").Value = "  -0.97%  "

# Row 26
$ws.Range("D26").Value = "1.787"
$ws.Range("E26").Value = "  -8.33%  "

# Row 27
$ws.Range("D27").Value = "17.99"
$ws.Range("E27").Value = "  -0.81%  "

# Row 28
$ws.Range("D28").Value = "5.088"
$ws.Range("E28").Value = "  +1.47%  "

# Row 29
$ws.Range("D29").Value = "114.18"
$ws.Range("E29").Value = "  +1.52%  "

# Row 30
$ws.Range("D30").Value = "1.837"
$ws.Range("E30").Value = "  +11.46%  "

# Row 31
$ws.Range("D31").Value = "0.08851"
$ws.Range("E31").Value = "  +0.01%  "

# Row 32
$ws.Range("D32").Value = "0.7276"
$ws.Range("E32").Value = "  +1.24%  "

# Row 33
$ws.Range("D33").Value = "1.122"
$ws.Range("E33").Value = "  +4.24%  "

# Row 34
$ws.Range("E34").Value = "  -0.41%  "

# Row 35
$ws.Range("E35").Value = "  -0.82%  "

# Row 36
$ws.Range("D36").Value = "2.732"
$ws.Range("E36").Value = "  -4.67%  "

# Row 37
$ws.Range("D37").Value = "1.081"
$ws.Range("E37").Value = "  +0.51%  "

# Row 38
$ws.Range("D38").Value = "0.05109"
$ws.Range("E38").Value = "  +0.47%  "

# Row 39
$ws.Range("D39").Value = "0.01887"
$ws.Range("E39").Value = "  -0.11%  "

# Row 40
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "0.4916"
$ws.Range("E40").Value = "  -0.16%  "

# Row 41
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").Value = "0.1608"
$ws.Range("E41").Value = "  -0.23%  "

# Row 42
$ws.Range("D42").Value = "2.599"
$ws.Range("E42").Value = "  +0.90%  "

# Row 43
$ws.Range("E43").Value = "  +2.51%  "

# Row 44
$ws.Range("D44").Value = "8.035"
$ws.Range("E44").Value = "  -0.15%  "

# Row 45
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "104.67"
$ws.Range("E45").Value = "  +0.49%  "

# Row 46
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "1.001"
$ws.Range("E46").Value = "  -0.67%  "

# Row 47
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "10.23"
$ws.Range("E47").Value = "  +0.83%  "

# Row 48
$ws.Range("D48").Value = "1.627"
$ws.Range("E48").Value = "  +2.72%  "

# Row 49
$ws.Range("D49").Value = "0.06171"
$ws.Range("E49").Value = "  -2.18%  "

# Row 50
$ws.Range("D50").Value = "0.4445"
$ws.Range("E50").Value = "  -1.57%  "

# Row 51
$ws.Range("D51").Value = "1.719"
$ws.Range("E51").Value = "  +4.19%  "

